$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.026.08'
$ws.Range("E2").Value = '  +0.23%  '

$ws.Range("D3").Value = '1.862.94'
$ws.Range("E3").Value = '  -0.27%  '

$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  +0.31%  '

$ws.Range("D5").Value = '''311.97'
$ws.Range("E5").Value = '  -0.16%  '

$ws.Range("D6").Value = '''1.003'

$ws.Range("E7").Value = '  +1.21%  '

$ws.Range("D8").Value = '''0.3828'
$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").Value = '''0.08300'
$ws.Range("E9").Value = '  -7.03%  '

$ws.Range("E10").Value = '  -0.19%  '

$ws.Range("D11").Value = '''41.51'
$ws.Range("E11").Value = '  +0.06%  '

$ws.Range("D12").Value = '''6.224'
$ws.Range("E12").Value = '  -2.25%  '

$ws.Range("D13").Value = '''20.61'
$ws.Range("E13").Value = '  -0.29%  '

$ws.Range("D14").Value = '1.856.79'
$ws.Range("E14").Value = '  -1.03%  '

$ws.Range("D15").Value = '''7.216'
$ws.Range("E15").Value = '  -0.26%  '

$ws.Range("D16").Value = '''1.004'
$ws.Range("E16").Value = '  +0.22%  '

$ws.Range("D17").Value = '''0.00001098'
$ws.Range("E17").Value = '  -0.10%  '

$ws.Range("E18").Value = '  -0.14%  '

$ws.Range("D19").Value = '''0.06632'
$ws.Range("E19").Value = '  -0.32%  '

$ws.Range("D20").Value = '''17.70'
$ws.Range("E20").Value = '  -2.47%  '

$ws.Range("E21").Value = '  +0.21%  '

$ws.Range("D22").Value = '''6.037'
$ws.Range("E22").Value = '  -1.17%  '

$ws.Range("D23").Value = '28.052.76'
$ws.Range("E23").Value = '  +0.26%  '

$ws.Range("D24").Value = '''11.09'
$ws.Range("E24").Value = '  -3.43%  '

$ws.Range("D25").Value = '''2.234'
$ws.Range("E25").Value = '  -1.66%  '

$ws.Range("D26").Value = '''2.547'
$ws.Range("E26").Value = '  +2.36%  '

$ws.Range("D27").Value = '2.074.90'
$ws.Range("E27").Value = '  -0.36%  '

$ws.Range("D28").Value = '''158.16'
$ws.Range("E28").Value = '  +0.28%  '

$ws.Range("D29").Value = '''20.53'
$ws.Range("E29").Value = '  -0.62%  '

$ws.Range("D30").Value = '''124.73'
$ws.Range("E30").Value = '  -1.00%  '

$ws.Range("D31").Value = '''0.1055'
$ws.Range("E31").Value = '  -0.75%  '

$ws.Range("D32").Value = '''1.041'
$ws.Range("E32").Value = '  -1.06%  '

$ws.Range("D33").Value = '''5.834'
$ws.Range("E33").Value = '  +4.27%  '

$ws.Range("D34").Value = '''3.587'
$ws.Range("E34").Value = '  -0.50%  '

$ws.Range("D35").Value = '''9.454'
$ws.Range("E35").Value = '  -0.20%  '

$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '''0.02419'
$ws.Range("E36").Value = '  +0.67%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = '''0.06525'
$ws.Range("E37").Value = '  -0.55%  '

$ws.Range("D38").Value = '''0.2174'
$ws.Range("E38").Value = '  -0.28%  '

$ws.Range("D39").Value = '''1.206'
$ws.Range("E39").Value = '  +0.48%  '

$ws.Range("D40").Value = '''0.6462'

$ws.Range("D41").Value = '''1.227'
$ws.Range("E41").Value = '  -4.38%  '

$ws.Range("D42").Value = '''4.947'
$ws.Range("E42").Value = '  +1.21%  '

$ws.Range("D43").Value = '''11.21'
$ws.Range("E43").Value = '  -2.44%  '

$ws.Range("D44").Value = '''0.6099'
$ws.Range("E44").Value = '  +1.61%  '

$ws.Range("D45").Value = '''13.12'
$ws.Range("E45").Value = '  -0.50%  '

$ws.Range("D46").Value = '''1.285'
$ws.Range("E46").Value = '  +0.53%  '

$ws.Range("E47").Value = '  +0.27%  '

$ws.Range("D48").Value = '''2.016'
$ws.Range("E48").Value = '  +1.38%  '

$ws.Range("D49").Value = '''1.208'
$ws.Range("E49").Value = '  -1.18%  '

$ws.Range("D50").Value = '''120.33'
$ws.Range("E50").Value = '  -0.24%  '

$ws.Range("D51").Value = '''78.37'
$ws.Range("E51").Value = '  -0.96%  '
